$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 418
$ws.Range("F3").Value = 411
$ws.Range("F4").Value = 2681
$ws.Range("F5").Value = 1312
$ws.Range("F6").Value = 54
$ws.Range("F8").Value = 563
$ws.Range("F9").Value = 35
$ws.Range("F10").Value = 573
$ws.Range("F12").Value = 77
$ws.Range("F13").Value = 11131
$ws.Range("F14").Value = 6359
$ws.Range("F21").Value = 878
$ws.Range("F23").Value = 229
$ws.Range("F24").Value = 894
$ws.Range("F25").Value = 3605
$ws.Range("F28").Value = 491
$ws.Range("F29").Value = 153
$ws.Range("F31").Value = 260
$ws.Range("F32").Value = 280
$ws.Range("F33").Value = 4941
$ws.Range("F35").Value = 1203
$ws.Range("F36").Value = 193
$ws.Range("F37").Value = 368
$ws.Range("F38").Value = 158
$ws.Range("F39").Value = 522

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G6").Value = "不可售"
$ws.Range("F12").Value = 3639
$ws.Range("F13").Value = 82

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8940
$ws.Range("F3").Value = 470
$ws.Range("F4").Value = 1751

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 8940
$ws.Range("F3").Value = 470
$ws.Range("F4").Value = 1751
$ws.Range("F6").Value = 411
$ws.Range("F7").Value = 2682
$ws.Range("F10").Value = 1312
$ws.Range("F12").Value = 54
$ws.Range("F13").Value = 563
$ws.Range("F14").Value = 35
$ws.Range("F16").Value = 573
$ws.Range("F18").Value = 77
$ws.Range("F19").Value = 11131
$ws.Range("F20").Value = 3639
$ws.Range("F21").Value = 6359
$ws.Range("F22").Value = 82
$ws.Range("F29").Value = 878
$ws.Range("F31").Value = 229
$ws.Range("F32").Value = 894
$ws.Range("F33").Value = 3605
$ws.Range("F35").Value = 153
$ws.Range("F37").Value = 260
$ws.Range("F40").Value = 280
$ws.Range("F41").Value = 4941
$ws.Range("F43").Value = 1203
$ws.Range("F44").Value = 193
$ws.Range("F45").Value = 158
$ws.Range("F46").Value = 522
